# Fixed failure scripts - Vacc & eCare. Updated Docker related files
#
# This script:
#  1. Removes the "Start Time" / "End Time" / "Duration" columns (D:F) from
#     every worksheet (eCare_POC, Negative_Scenario, OnDemand_POC,
#     Vaccinations_POC) - those columns are no longer used by the report.
#  2. Refreshes the randomly-generated test data (Firstname/Lastname values)
#     on the eCare_POC and Vaccinations_POC sheets.
#  3. Adds a new test step on the Vaccinations_POC sheet documenting
#     validation of the Vaccine History section after adding a vaccination.

$wb = $excel.ActiveWorkbook

$sheetNames = @("eCare_POC", "Negative_Scenario", "OnDemand_POC", "Vaccinations_POC")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns("D:F").Delete()
}

# --- eCare_POC: refresh the patient first/last name test data ------------
$wsECare = $wb.Worksheets.Item("eCare_POC")
$wsECare.Range("B7").Value = " Firstname:TestAutomation6840"
$wsECare.Range("B8").Value = " Lastname:User92986"

# --- Vaccinations_POC: refresh patient first/last name + add new step ----
$wsVacc = $wb.Worksheets.Item("Vaccinations_POC")
$wsVacc.Range("B7").Value = " Firstname:TestAutomation48623"
$wsVacc.Range("B8").Value = " Lastname:User4558"

# Insert a new row 21 for the additional validation step, pushing the
# trailing blank row down to row 22.
$wsVacc.Rows("21:21").Insert()
$wsVacc.Range("A21").Value = "Step 13"
$wsVacc.Range("B21").Value = " Validated Vaccine History Section for added vaccination "
$wsVacc.Range("C21").Value = " Passed"

$wb.Save()
